$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# Re-order the three RACE/ETHNICITY rows (12-14): native american / asian / pacific islander
# become pacific islander / native american / asian, along with their associated
# "Predominantly CSA" (D) and "Mainly OSA" (F) counts swapping to match.
$ws.Range("A12").Value = "pacific islander"

$ws.Range("A13").Value = "native american"
$ws.Range("D13").Value = "1/62 (1.6%)"
$ws.Range("F13").Value = "1/171 (0.6%)"

$ws.Range("A14").Value = "asian"
$ws.Range("D14").Value = "0/62 (0.0%)"
$ws.Range("F14").Value = "2/171 (1.2%)"

# Updated OSA severity figures (Sept 2021 submission) for "mild" (row 35) and
# "none" (row 36) severities.
$ws.Range("B35").Value = "93/510 (18.2%)"
$ws.Range("E35").Value = "55/266 (20.7%)"
$ws.Range("F35").Value = "30/171 (17.5%)"

$ws.Range("B36").Value = "0/510 (0.0%)"
$ws.Range("E36").Value = "0/266 (0.0%)"
$ws.Range("F36").Value = "0/171 (0.0%)"
